# Initial thermometer calibration setup
# Insert a new localization row ("thermometer_calibration" /
# "Thermometer Calibration") right above the existing "food_beef" row
# (i.e. as the new row 15), pushing every row below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Insert a blank row at position 15 - everything from the old row 15
# downward (food_beef ... washProduce_item_7) shifts down to make room.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row with the new key/value pair.
$ws.Range("A15").Value = "thermometer_calibration"
$ws.Range("B15").Value = "Thermometer Calibration"

# Match the author's recorded view state: scrolled down a bit with the
# new row selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("A15").Select()
